$d = $word.ActiveDocument

function Split-LabelField($label, $placeholder) {
    # $label already contains the literal "Prefix: Value" text currently in the
    # document; this finds it, trims it back to "Prefix: " and appends a
    # newly-split bold run containing the placeholder token.
    $prefix = $label.Prefix
    $old = $label.Old

    $rng = $d.Content
    $rng.Find.Execute($old) | Out-Null
    $rng.Text = $prefix
    $rng.Collapse(0)
    $rng.InsertAfter($placeholder)
    # Toggling Bold off/on forces the newly inserted text to live in its own
    # run instead of being silently re-absorbed into the preceding run.
    $rng.Bold = $false
    $rng.Bold = $true
}

# 1) Nacionalidade: BRASILEIRO(A)  ->  Nacionalidade: {{nacionalidade}}
Split-LabelField @{ Old = "Nacionalidade: BRASILEIRO(A)"; Prefix = "Nacionalidade: " } "{{nacionalidade}}"

# 2) Estado Civil value placeholder {{data}} -> {{estadoCivil}}  (do this before
#    we introduce a new {{data}} token further down, so Find stays unambiguous).
#    This token is already its own run; a plain text swap would get silently
#    re-absorbed into the previous ("Estado Civil: ") run, so use the same
#    bold off/on trick to keep it split out.
$rng = $d.Content
$rng.Find.Execute("{{data}}") | Out-Null
$rng.Text = "{{estadoCivil}}"
$rng.Bold = $false
$rng.Bold = $true

# 3) Profissão: PEDREIRO  ->  Profissão: {{profissão}}
Split-LabelField @{ Old = "Profissão: PEDREIRO"; Prefix = "Profissão: " } "{{profissão}}"

# 4) FONE: ...  ->  FONE: {{fone}}
Split-LabelField @{ Old = "FONE: 47-99733-4601 OU 47 9 9976-8577"; Prefix = "FONE: " } "{{fone}}"

# 5) CPF number -> {{cpf}} (already its own run after "Nº CPF: "; keep it split)
$rng = $d.Content
$rng.Find.Execute("076.336.989-69") | Out-Null
$rng.Text = "{{cpf}}"
$rng.Bold = $false
$rng.Bold = $true

# 6) Nº RG: 55519378  ->  Nº RG: {{rg}}
Split-LabelField @{ Old = "Nº RG: 55519378"; Prefix = "Nº RG: " } "{{rg}}"

# 7) Endereço: ...  ->  Endereço: {{endereco}}
Split-LabelField @{ Old = "Endereço: RUA OGIDIO DA SILVA N32 BLOCO 12 AP 13"; Prefix = "Endereço: " } "{{endereco}}"

# 8) Bairro: COLONINHA  ->  Bairro: {{bairro}}
Split-LabelField @{ Old = "Bairro: COLONINHA"; Prefix = "Bairro: " } "{{bairro}}"

# 9) CEP: 89.110-260  ->  CEP: {{cep}}
Split-LabelField @{ Old = "CEP: 89.110-260"; Prefix = "CEP: " } "{{cep}}"

# 10) CIDADE: GASPAR  ->  CIDADE: {{cidade}}
Split-LabelField @{ Old = "CIDADE: GASPAR"; Prefix = "CIDADE: " } "{{cidade}}"

# 11) Estado: SANTA CATARINA  ->  Estado: {{estado}}
Split-LabelField @{ Old = "Estado: SANTA CATARINA"; Prefix = "Estado: " } "{{estado}}"

# 12) Replace the "Jaraguá do Sul - SC, <TIME field>." construct with
#     "Jaraguá do Sul - SC,{{data}}."
$i = 0
foreach ($f in $d.Fields) {
    $i = $i + 1
    if ($i -eq 2) {
        $f.Delete()
    }
}

$rng = $d.Content
$rng.Find.Execute("Jaraguá do Sul - SC, ") | Out-Null
$rng.Text = "Jaraguá do Sul - SC,"
$rng.Collapse(0)
$rng.InsertAfter("{{data}}")
# Here the placeholder must NOT be bold, so toggle the opposite direction.
$rng.Bold = $true
$rng.Bold = $false

Write-Host "done"
